$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Power Consumption (mW)" column ---
$ws.Range("C1").Value = "Power Consumption (mW)"

$ws.Range("C2").Formula = "=3.3*B2"
$ws.Range("C3:C6").Formula = "=3.3*B3"

# C2:C5 get the same thin box border the rest of the table uses
$ws.Range("C2:C5").Borders.LineStyle = 1
$ws.Range("C2:C5").Borders.Weight = 2

# C5 sits on the rule above the Total row, so its bottom edge is doubled
$ws.Range("C5").Borders.Item(9).LineStyle = -4119
$ws.Range("C5").Borders.Item(9).Weight = -4138

# Column B no longer sits on the right edge of the table (column C does),
# so drop its right border for the data rows - the header keeps its box.
$ws.Range("B2:B5").Borders.Item(10).LineStyle = -4142

# Re-select so the active cell matches the newly-added column
$ws.Range("C1:C6").Select()
